$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 53

$ws.Cells.Item($row, 1).Value = 1
$ws.Cells.Item($row, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item($row, 3).Value = "Arica y Parinacota"

# Date: 2022-12-16 (serial 44911) - reuse the same number format as the cell above
# so the engine doesn't register a new (unused) style.
$ws.Cells.Item($row, 4).NumberFormat = $ws.Cells.Item($row - 1, 4).NumberFormat
$ws.Cells.Item($row, 4).Value = (Get-Date -Year 2022 -Month 12 -Day 16 -Hour 0 -Minute 0 -Second 0)

$ws.Cells.Item($row, 5).Value = 15
$ws.Cells.Item($row, 6).Value = 100112028
$ws.Cells.Item($row, 7).Value = "Sandia"
$ws.Cells.Item($row, 8).Value = "Sin especificar"
$ws.Cells.Item($row, 9).Value = "Primera"
$ws.Cells.Item($row, 10).Value = 500
$ws.Cells.Item($row, 11).Value = 350
$ws.Cells.Item($row, 12).Value = 380
$ws.Cells.Item($row, 13).Value = 365
$ws.Cells.Item($row, 14).Value = '$/kilo (volumen en unidades)'
$ws.Cells.Item($row, 15).Value = "Perú"
$ws.Cells.Item($row, 16).Value = 365
$ws.Cells.Item($row, 17).Value = 1
$ws.Cells.Item($row, 18).Value = "Hortaliza"
